$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-17 Tuesday" "2024-09-18 Wednesday"

Replace-Text "761×5=" "114×3="
Replace-Text "230×3=" "685×6="
Replace-Text "367×5=" "106×9="
Replace-Text "844×9=" "530×7="
Replace-Text "418×2=" "744×6="

Replace-Text "181×5=" "661×5="
Replace-Text "261×3=" "680×8="
Replace-Text "423×6=" "494×5="
Replace-Text "721×3=" "127×2="
Replace-Text "935×2=" "315×6="

Replace-Text "698×2=" "278×3="
Replace-Text "484×8=" "800×6="
Replace-Text "513×2=" "736×7="
Replace-Text "814×4=" "197×2="
Replace-Text "620×9=" "447×7="

Replace-Text "577×6=" "751×3="
Replace-Text "944×2=" "799×5="
Replace-Text "609×5=" "675×2="
Replace-Text "144×3=" "782×8="
Replace-Text "199×3=" "455×7="

Replace-Text "612×5=" "503×8="
Replace-Text "138×6=" "436×3="
Replace-Text "891×7=" "365×8="
Replace-Text "759×9=" "760×8="
Replace-Text "239×6=" "899×3="
